$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 52.75
$ws.Range("I5").Value = 52.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 52.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 62.25

$ws.Range("H15").Value = 1249.3625
$ws.Range("I15").Value = 1249.3625
$ws.Range("K15").Value = 3748.0875
$ws.Range("M15").Value = -3579.0875

$ws.Range("H18").Value = 715.2105
$ws.Range("I18").Value = 588.2778
$ws.Range("K18").Value = 588.2778
$ws.Range("M18").Value = -304.2778

$ws.Range("H40").Value = 1680.4615
$ws.Range("J40").Value = 1985.7142
$ws.Range("L40").Value = 1985.7142
$ws.Range("N40").Value = -2335.7142

$ws.Range("H41").Value = 1924.5
$ws.Range("I41").Value = 3020
$ws.Range("J41").Value = 1142
$ws.Range("K41").Value = 3020
$ws.Range("L41").Value = 1142
$ws.Range("M41").Value = -2580
$ws.Range("N41").Value = -2022

$ws.Range("H111").Value = 6988.8096
$ws.Range("I111").Value = 13895.875
$ws.Range("J111").Value = 2738.3076
$ws.Range("K111").Value = 41687.625
$ws.Range("L111").Value = 8214.9228
$ws.Range("M111").Value = -38620.625
$ws.Range("N111").Value = -14348.9228

$ws.Range("H116").Value = 2081.8948
$ws.Range("I116").Value = 1811.25
$ws.Range("J116").Value = 2278.7273
$ws.Range("K116").Value = 1811.25
$ws.Range("L116").Value = 2278.7273
$ws.Range("M116").Value = 1630.75
$ws.Range("N116").Value = -9162.7273

$ws.Range("H132").Value = 5439701.5
$ws.Range("I132").Value = 5686862
$ws.Range("J132").Value = 2169.5
$ws.Range("K132").Value = 17060586
$ws.Range("L132").Value = 6508.5
$ws.Range("M132").Value = -17058056
$ws.Range("N132").Value = -11568.5

$ws.Range("H137").Value = 1445.1852
$ws.Range("I137").Value = 1191.4348
$ws.Range("J137").Value = 1633.4517
$ws.Range("K137").Value = 3574.3044
$ws.Range("L137").Value = 4900.355100000001
$ws.Range("M137").Value = -1024.3044
$ws.Range("N137").Value = -10000.3551

$ws.Range("H138").Value = 5016.698
$ws.Range("I138").Value = 3987.2222
$ws.Range("J138").Value = 5227.273
$ws.Range("K138").Value = 11961.6666
$ws.Range("L138").Value = 15681.819
$ws.Range("M138").Value = -6821.6666
$ws.Range("N138").Value = -25961.819


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 198.5
$ws.Range("I5").Value = 178.2
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 178.2
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -66.19999999999999
$ws.Range("N5").Value = -524

$ws.Range("H32").Value = 18867.959
$ws.Range("I32").Value = 3293.5625
$ws.Range("J32").Value = 92159.234
$ws.Range("K32").Value = 3293.5625
$ws.Range("L32").Value = 92159.234
$ws.Range("M32").Value = -3006.5625
$ws.Range("N32").Value = -92733.234

$ws.Range("H41").Value = 3252
$ws.Range("I41").Value = 3252
$ws.Range("K41").Value = 3252
$ws.Range("M41").Value = -2838

$ws.Range("H45").Value = 2763.7273
$ws.Range("I45").Value = 2488.7368
$ws.Range("J45").Value = 3136.9285
$ws.Range("K45").Value = 2488.7368
$ws.Range("L45").Value = 3136.9285
$ws.Range("M45").Value = -2111.7368
$ws.Range("N45").Value = -3890.9285

$ws.Range("H61").Value = 1666.5883
$ws.Range("I61").Value = 834.5294
$ws.Range("J61").Value = 2082.6177
$ws.Range("K61").Value = 834.5294
$ws.Range("L61").Value = 2082.6177
$ws.Range("M61").Value = -622.5294
$ws.Range("N61").Value = -2506.6177

$ws.Range("H63").Value = 2492.9412
$ws.Range("I63").Value = 2169.2307
$ws.Range("K63").Value = 2169.2307
$ws.Range("M63").Value = -1483.2307

$ws.Range("H66").Value = 2492.9412
$ws.Range("I66").Value = 2169.2307
$ws.Range("K66").Value = 10846.1535
$ws.Range("M66").Value = -7414.1535

$ws.Range("H74").Value = 2979.3447
$ws.Range("I74").Value = 2349.6667
$ws.Range("K74").Value = 2349.6667
$ws.Range("M74").Value = -1475.6667

$ws.Range("H77").Value = 2979.3447
$ws.Range("I77").Value = 2349.6667
$ws.Range("K77").Value = 11748.3335
$ws.Range("M77").Value = -7380.333500000001

$ws.Range("H110").Value = 22774602
$ws.Range("I110").Value = 38539690
$ws.Range("K110").Value = 38539690
$ws.Range("M110").Value = -38537645

$ws.Range("H124").Value = 27425.428
$ws.Range("J124").Value = 27425.428
$ws.Range("L124").Value = 27425.428
$ws.Range("N124").Value = -37245.428

$ws.Range("H125").Value = 40857.5
$ws.Range("J125").Value = 40857.5
$ws.Range("L125").Value = 40857.5
$ws.Range("N125").Value = -50697.5

$ws.Range("H136").Value = 1666.5883
$ws.Range("I136").Value = 834.5294
$ws.Range("J136").Value = 2082.6177
$ws.Range("K136").Value = 2503.5882
$ws.Range("L136").Value = 6247.853099999999
$ws.Range("M136").Value = 46.41179999999986
$ws.Range("N136").Value = -11347.8531


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 198.5
$ws.Range("I4").Value = 178.2
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 178.2
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -63.19999999999999
$ws.Range("N4").Value = -530

$ws.Range("H105").Value = 184243.55
$ws.Range("I105").Value = 127183.5
$ws.Range("J105").Value = 336403.66
$ws.Range("K105").Value = 127183.5
$ws.Range("L105").Value = 336403.66
$ws.Range("M105").Value = -125436.5
$ws.Range("N105").Value = -339897.66

$ws.Range("H124").Value = 48000
$ws.Range("J124").Value = 48000
$ws.Range("L124").Value = 48000
$ws.Range("N124").Value = -57820


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96.78570999999999
$ws.Range("I7").Value = 46.2
$ws.Range("J7").Value = 223.25
$ws.Range("K7").Value = 46.2
$ws.Range("L7").Value = 223.25
$ws.Range("M7").Value = 66.8
$ws.Range("N7").Value = -449.25

$ws.Range("H22").Value = 704.75
$ws.Range("I22").Value = 327.6
$ws.Range("K22").Value = 327.6
$ws.Range("M22").Value = 22.39999999999998

$ws.Range("H31").Value = 14105.827
$ws.Range("I31").Value = 27335.29
$ws.Range("J31").Value = 2414.6743
$ws.Range("K31").Value = 27335.29
$ws.Range("L31").Value = 2414.6743
$ws.Range("M31").Value = -27040.29
$ws.Range("N31").Value = -3004.6743

$ws.Range("H34").Value = 14105.827
$ws.Range("I34").Value = 27335.29
$ws.Range("J34").Value = 2414.6743
$ws.Range("K34").Value = 27335.29
$ws.Range("L34").Value = 2414.6743
$ws.Range("M34").Value = -27133.29
$ws.Range("N34").Value = -2818.6743

$ws.Range("H51").Value = 7971.6
$ws.Range("J51").Value = 7971.6
$ws.Range("L51").Value = 7971.6
$ws.Range("N51").Value = -9443.6

$ws.Range("H61").Value = 7971.6
$ws.Range("J61").Value = 7971.6
$ws.Range("L61").Value = 7971.6
$ws.Range("N61").Value = -8667.6

$ws.Range("H74").Value = 21304.416
$ws.Range("J74").Value = 21304.416
$ws.Range("L74").Value = 21304.416
$ws.Range("N74").Value = -23052.416

$ws.Range("H77").Value = 21304.416
$ws.Range("J77").Value = 21304.416
$ws.Range("L77").Value = 63913.24800000001
$ws.Range("N77").Value = -72649.24800000001

$ws.Range("H99").Value = 8841.549999999999
$ws.Range("I99").Value = 4037.8
$ws.Range("J99").Value = 13645.3
$ws.Range("K99").Value = 4037.8
$ws.Range("L99").Value = 13645.3
$ws.Range("M99").Value = -2539.8
$ws.Range("N99").Value = -16641.3

$ws.Range("H124").Value = 33420
$ws.Range("J124").Value = 33420
$ws.Range("L124").Value = 33420
$ws.Range("N124").Value = -38330

$ws.Range("H126").Value = 8841.549999999999
$ws.Range("I126").Value = 4037.8
$ws.Range("J126").Value = 13645.3
$ws.Range("K126").Value = 12113.4
$ws.Range("L126").Value = 40935.89999999999
$ws.Range("M126").Value = -9643.400000000001
$ws.Range("N126").Value = -45875.89999999999

$ws.Range("H132").Value = 2799.8333
$ws.Range("I132").Value = 2765.1304
$ws.Range("K132").Value = 8295.3912
$ws.Range("M132").Value = -5765.3912

$ws.Range("H134").Value = 2373.45
$ws.Range("I134").Value = 1794.3846
$ws.Range("J134").Value = 3448.8572
$ws.Range("K134").Value = 5383.1538
$ws.Range("L134").Value = 10346.5716
$ws.Range("M134").Value = -2848.1538
$ws.Range("N134").Value = -15416.5716


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1997.1351
$ws.Range("I68").Value = 1311.6562
$ws.Range("K68").Value = 3934.9686
$ws.Range("M68").Value = -3123.9686

$ws.Range("H71").Value = 1997.1351
$ws.Range("I71").Value = 1311.6562
$ws.Range("K71").Value = 11804.9058
$ws.Range("M71").Value = -7748.905799999999

$ws.Range("H107").Value = 290665.1
$ws.Range("I107").Value = 647.3077
$ws.Range("J107").Value = 633413.4399999999
$ws.Range("K107").Value = 1941.9231
$ws.Range("L107").Value = 1900240.32
$ws.Range("M107").Value = -21.92309999999998
$ws.Range("N107").Value = -1904080.32

$ws.Range("H131").Value = 1131301.6
$ws.Range("J131").Value = 1283433.8
$ws.Range("L131").Value = 3850301.4
$ws.Range("N131").Value = -3860381.4

$ws.Range("H140").Value = 1535.409
$ws.Range("J140").Value = 2549.8333
$ws.Range("L140").Value = 7649.499899999999
$ws.Range("N140").Value = -18009.4999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 66668670

$ws.Range("H132").Value = 1814
$ws.Range("I132").Value = 1310.0526
$ws.Range("J132").Value = 2771.5
$ws.Range("K132").Value = 3930.1578
$ws.Range("L132").Value = 8314.5
$ws.Range("M132").Value = -1400.1578
$ws.Range("N132").Value = -13374.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5104.6313
$ws.Range("I132").Value = 6135.636
$ws.Range("J132").Value = 3687
$ws.Range("K132").Value = 18406.908
$ws.Range("L132").Value = 11061
$ws.Range("M132").Value = -15876.908
$ws.Range("N132").Value = -16121

$ws.Range("H136").Value = 1830.7142
$ws.Range("I136").Value = 1450.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4351.200000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1801.200000000001
$ws.Range("N136").Value = -20100

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0


# Row-level cell removals (diff removes these cells entirely)
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("N5").ClearContents()

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("N139").ClearContents()
